# Master_CON_GS_CO2.xlsx - apply Gasera CO2 -> CCO2 (carbon-mass) transformation
# Commit: "CCH4, NN2O and CCO2 Gasera transformations in Master_GHG_2023"
#
# The Gasera CO2 flux columns (C: avg_Gasera_CO2_flux_mgm2h,
# D: avg_Gasera_CO2_flux_mgm2h_cor) are renamed to the "CCO2" (carbon-mass)
# variant and their values are converted from mg CO2 / m2 / h to the
# corresponding carbon-mass units by dividing by (44/12)^2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames -------------------------------------------------------
$ws.Range("C1").Value = "avg_Gasera_CCO2_flux_mgm2h"
$ws.Range("D1").Value = "avg_Gasera_CCO2_flux_mgm2h_cor"

# --- Converted data values -------------------------------------------------
$newC = @{
    2  = -154.6954613750029
    3  = -132.7707040094541
    4  = -130.7769647958405
    5  = -146.3263235306765
    6  = -129.6483883732642
    7  = -51.51541126762701
    8  = -89.04079704369862
    9  = -53.60230070600097
    10 = -92.48860176531113
    11 = -15.34460319066201
    12 = -70.79988975687567
    13 = -25.4462588262617
    14 = -88.19401974602674
    15 = 8.160756652706471
    16 = -76.0331491692946
    17 = -86.05664164143167
    18 = -69.59577951479933
}

$newD = @{
    2  = -154.6954613750029
    3  = -132.7707040094541
    4  = -130.7769647958405
    5  = -146.3263235306765
    6  = -129.6483883732642
    7  = -51.51541126762701
    8  = -89.04079704369862
    9  = -53.60230070600097
    10 = -92.48860176531113
    11 = -15.34460319066201
    12 = -70.79988975687567
    13 = -25.4462588262617
    14 = -88.19401974602674
    # row 15 already held a 0 override (pre-correction) and stays 0
    15 = 0
    16 = -76.0331491692946
    17 = -86.05664164143167
    18 = -69.59577951479933
}

foreach ($row in $newC.Keys) {
    $ws.Cells.Item($row, 3).Value = $newC[$row]
}

foreach ($row in $newD.Keys) {
    $ws.Cells.Item($row, 4).Value = $newD[$row]
}
